# Update "想去人数" (want-to-go count) figures in column F on the
# "展览" and "全部类型" sheets to the newly scraped values.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 71
$ws1.Range("F3").Value = 11757
$ws1.Range("F5").Value = 342
$ws1.Range("F7").Value = 11694
$ws1.Range("F8").Value = 484
$ws1.Range("F10").Value = 95
$ws1.Range("F11").Value = 29
$ws1.Range("F12").Value = 1767
$ws1.Range("F13").Value = 5791
$ws1.Range("F14").Value = 119
$ws1.Range("F15").Value = 3523

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 71
$ws4.Range("F5").Value = 11757
$ws4.Range("F7").Value = 342
$ws4.Range("F9").Value = 11694
$ws4.Range("F10").Value = 484
$ws4.Range("F12").Value = 95
$ws4.Range("F13").Value = 29
$ws4.Range("F14").Value = 1767
$ws4.Range("F16").Value = 5791
$ws4.Range("F17").Value = 119
$ws4.Range("F18").Value = 3523
